$ws = $excel.ActiveWorkbook.ActiveSheet

# Rows 61-119 correspond to Bono-A2 (ISIN ES0377955010); NUM_BONOS (col G)
# changed from 100 to 589, with dependent totals in columns M and N rescaled.
$data = New-Object 'object[][]' 59
$data[0] = @(61, 589, 4565869.1, 0)
$data[1] = @(62, 589, 4243214.9, 24844.02)
$data[2] = @(63, 589, 3891570.12, 23088.8)
$data[3] = @(64, 589, 3554632.67, 20715.13)
$data[4] = @(65, 589, 3233309.61, 19130.72)
$data[5] = @(66, 589, 2928213.5, 17593.43)
$data[6] = @(67, 589, 2636664.39, 15932.45)
$data[7] = @(68, 589, 2358833.09, 14035.87)
$data[8] = @(69, 589, 2091480.1, 12692.95)
$data[9] = @(70, 589, 1834788.01, 11379.48)
$data[10] = @(71, 589, 1587437.46, 9983.549999999999)
$data[11] = @(72, 589, 1350241.27, 8540.5)
$data[12] = @(73, 589, 1122015.55, 7268.26)
$data[13] = @(74, 589, 902948.78, 6237.51)
$data[14] = @(75, 589, 692616.88, 4859.25)
$data[15] = @(76, 589, 490536.87, 3728.37)
$data[16] = @(77, 589, 296420.14, 2638.72)
$data[17] = @(78, 589, 109972.19, 1596.19)
$data[18] = @(79, 589, 0, 589)
$data[19] = @(80, 589, 4565869.1, 0)
$data[20] = @(81, 589, 4267605.39, 24844.02)
$data[21] = @(82, 589, 3985003.19, 23218.38)
$data[22] = @(83, 589, 3711165.31, 21209.89)
$data[23] = @(84, 589, 3447387.55, 19972.99)
$data[24] = @(85, 589, 3194682.99, 18753.76)
$data[25] = @(86, 589, 2950636.73, 17381.39)
$data[26] = @(87, 589, 2715772.98, 15702.74)
$data[27] = @(88, 589, 2486958.26, 14613.09)
$data[28] = @(89, 589, 2264634.32, 13529.33)
$data[29] = @(90, 589, 2047629.05, 12321.88)
$data[30] = @(91, 589, 1837055.66, 11020.19)
$data[31] = @(92, 589, 1631865.73, 9889.309999999999)
$data[32] = @(93, 589, 1432483.34, 9070.6)
$data[33] = @(94, 589, 1238661.11, 7710.01)
$data[34] = @(95, 589, 1050080.98, 6667.48)
$data[35] = @(96, 589, 866607.48, 5648.51)
$data[36] = @(97, 589, 688116.92, 4664.88)
$data[37] = @(98, 589, 514532.73, 3704.81)
$data[38] = @(99, 589, 345619.31, 2768.3)
$data[39] = @(100, 589, 181465.01, 1861.24)
$data[40] = @(101, 589, 22246.53, 977.74)
$data[41] = @(102, 589, -0, 123.69)
$data[42] = @(103, 589, 4565869.1, 0)
$data[43] = @(104, 589, 4222140.48, 24844.02)
$data[44] = @(105, 589, 3811572.14, 22971)
$data[45] = @(106, 589, 3421819.06, 20285.16)
$data[46] = @(107, 589, 3053322.99, 18412.14)
$data[47] = @(108, 589, 2706213.51, 16609.8)
$data[48] = @(109, 589, 2377463.16, 14725)
$data[49] = @(110, 589, 2066842.23, 12651.72)
$data[50] = @(111, 589, 1770910.96, 11120.32)
$data[51] = @(112, 589, 1489527.99, 9636.040000000001)
$data[52] = @(113, 589, 1221173.7, 8104.64)
$data[53] = @(114, 589, 966313.4, 6573.24)
$data[54] = @(115, 589, 723586.5, 5200.87)
$data[55] = @(116, 589, 492910.54, 4022.87)
$data[56] = @(117, 589, 273643.51, 2650.5)
$data[57] = @(118, 589, 65119.84, 1472.5)
$data[58] = @(119, 589, -0, 347.51)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 7).Value = $item[1]   # G: NUM_BONOS
    $ws.Cells.Item($r, 13).Value = $item[2]  # M: TT1
    $ws.Cells.Item($r, 14).Value = $item[3]  # N: TT2
}